$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 0.04431377938692467
$ws.Cells.Item(2, 8).Value = 4.813587686683594
$ws.Cells.Item(2, 9).Value = -9.189826254501501

$ws.Cells.Item(3, 7).Value = 0.04525132216480877
$ws.Cells.Item(3, 8).Value = -7.082140266838807

$ws.Cells.Item(4, 7).Value = -0.004013984100506916
$ws.Cells.Item(4, 8).Value = -310.8300579558673

$ws.Cells.Item(5, 7).Value = 0.02528245319252126
$ws.Cells.Item(5, 8).Value = 565.9672329208537

$ws.Cells.Item(6, 7).Value = 0.01836266350497041
$ws.Cells.Item(6, 8).Value = -47.03522859908454

$ws.Cells.Item(7, 7).Value = 0.04048167953056499
$ws.Cells.Item(7, 8).Value = -23.89080394777568

$ws.Cells.Item(8, 7).Value = -0.03874419896777712
$ws.Cells.Item(8, 8).Value = -105.8586796609504

$ws.Cells.Item(9, 7).Value = -0.008741720791396387
$ws.Cells.Item(9, 8).Value = 59.39353696687668

$ws.Cells.Item(10, 7).Value = -0.09693825342544379
$ws.Cells.Item(10, 8).Value = -33.33761398313414

$ws.Cells.Item(11, 7).Value = -0.07398412036013222
$ws.Cells.Item(11, 8).Value = 19.59212815136385

$ws.Cells.Item(12, 7).Value = -0.2225324303338621
$ws.Cells.Item(12, 8).Value = 8.967804220329887

$ws.Cells.Item(13, 7).Value = -0.3324163756192516
$ws.Cells.Item(13, 8).Value = -20.96083823775607

$ws.Cells.Item(14, 7).Value = -0.02604287913974952
$ws.Cells.Item(14, 8).Value = 29.7982499302573

$ws.Cells.Item(15, 7).Value = 0.008403580477210248
$ws.Cells.Item(15, 8).Value = 124.1681103211719

$ws.Cells.Item(16, 7).Value = 0.1182794885734315
$ws.Cells.Item(16, 8).Value = -5.609129793604348

$ws.Cells.Item(17, 7).Value = 0.1269171858770419
$ws.Cells.Item(17, 8).Value = -9.506779757116151

$ws.Cells.Item(18, 7).Value = 0.1142783311011611
$ws.Cells.Item(18, 8).Value = -8.38143557566344

$ws.Cells.Item(19, 7).Value = 0.127761031392099
$ws.Cells.Item(19, 8).Value = -4.093276920888729

$ws.Cells.Item(20, 7).Value = 0.0251666992917923
$ws.Cells.Item(20, 8).Value = -26.70486755209642

$ws.Cells.Item(21, 7).Value = 0.06941043298523165
$ws.Cells.Item(21, 8).Value = 19.59096706885674

$ws.Cells.Item(22, 7).Value = -0.07401208302312082
$ws.Cells.Item(22, 8).Value = 7.308764466275205

$ws.Cells.Item(23, 7).Value = -0.07525245114147985
$ws.Cells.Item(23, 8).Value = -20.31453406654256

$ws.Cells.Item(24, 7).Value = 0.1099801663812027
$ws.Cells.Item(24, 8).Value = -6.886828429307662

$ws.Cells.Item(25, 7).Value = 0.1598671706431533
$ws.Cells.Item(25, 8).Value = 26.70735653795628

$ws.Cells.Item(26, 7).Value = 0.05044056829285298
$ws.Cells.Item(26, 8).Value = 1.481688569075823

$ws.Cells.Item(27, 7).Value = 0.0899068080799017
$ws.Cells.Item(27, 8).Value = 3.72677156575148

$ws.Cells.Item(28, 7).Value = -0.08401321649715547
$ws.Cells.Item(28, 8).Value = -32.10863964441598

$ws.Cells.Item(29, 7).Value = -0.09693873732094824
$ws.Cells.Item(29, 8).Value = -36.20305196176385

$ws.Cells.Item(30, 7).Value = 0.08061527957701278
$ws.Cells.Item(30, 8).Value = 26.5384553489482

$ws.Cells.Item(31, 7).Value = 0.05451001821571826
$ws.Cells.Item(31, 8).Value = -10.02046356563199

$ws.Cells.Item(32, 7).Value = 0.08349834395145604
$ws.Cells.Item(32, 8).Value = -15.02650275945227

$ws.Cells.Item(33, 7).Value = 0.1118461922945107
$ws.Cells.Item(33, 8).Value = 35.92585901937252

$ws.Cells.Item(34, 7).Value = -0.003021708923425605
$ws.Cells.Item(34, 8).Value = -111.5973282564692

$ws.Cells.Item(35, 7).Value = 0.002543253418372906
$ws.Cells.Item(35, 8).Value = 122.6871428614555

$ws.Cells.Item(36, 7).Value = -0.01282527203172271
$ws.Cells.Item(36, 8).Value = -2482.962334578323

$ws.Cells.Item(37, 7).Value = 0.001414960094942574
$ws.Cells.Item(37, 8).Value = 111.2707522229052

$ws.Cells.Item(38, 7).Value = 0.1128641674169546
$ws.Cells.Item(38, 8).Value = 5.227524548446874

$ws.Cells.Item(39, 7).Value = 0.1014787226493753
$ws.Cells.Item(39, 8).Value = 18.46302480226925

$ws.Cells.Item(40, 7).Value = 0.01203153017014169
$ws.Cells.Item(40, 8).Value = 305.0669354527938

$ws.Cells.Item(41, 7).Value = 0.02568259534296581
$ws.Cells.Item(41, 8).Value = 71.25444814297532

$ws.Cells.Item(42, 7).Value = 0.1169077193137913
$ws.Cells.Item(42, 8).Value = 15.82461378502358

$ws.Cells.Item(43, 7).Value = 0.1270727624769456
$ws.Cells.Item(43, 8).Value = 5.766605273878025

$ws.Cells.Item(44, 7).Value = 0.03781365682510281
$ws.Cells.Item(44, 8).Value = 5.957448863010271

$ws.Cells.Item(45, 7).Value = 0.03292585268647224
$ws.Cells.Item(45, 8).Value = 101.1331850092742

$ws.Cells.Item(46, 7).Value = 0.04648799332327846
$ws.Cells.Item(46, 8).Value = 28.28336925000987

$ws.Cells.Item(47, 7).Value = 0.06012295436664161
$ws.Cells.Item(47, 8).Value = 19.19623364338651

$ws.Cells.Item(48, 7).Value = 0.05422798928842538
$ws.Cells.Item(48, 8).Value = 26.75628479798734

$ws.Cells.Item(49, 7).Value = 0.0645193910199571
$ws.Cells.Item(49, 8).Value = -7.13440967921494

$ws.Cells.Item(50, 7).Value = 0.00157071224846759
$ws.Cells.Item(50, 8).Value = -90.90639182925533

$ws.Cells.Item(51, 7).Value = 0.0143775099409589
$ws.Cells.Item(51, 8).Value = -26.15447704044383

$ws.Cells.Item(52, 7).Value = -0.08852302347825509
$ws.Cells.Item(52, 8).Value = 14.48740915419254

$ws.Cells.Item(53, 7).Value = -0.08533288343904034
$ws.Cells.Item(53, 8).Value = 7.603404454407194

$ws.Cells.Item(54, 7).Value = 0.08205188786766718
$ws.Cells.Item(54, 8).Value = 12.21129907332394

$ws.Cells.Item(55, 7).Value = 0.08292468881285663
$ws.Cells.Item(55, 8).Value = 33.85386667813378

$ws.Cells.Item(56, 7).Value = 0.05528708507196726
$ws.Cells.Item(56, 8).Value = 58.0117907002854

$ws.Cells.Item(57, 7).Value = 0.016626581560588
$ws.Cells.Item(57, 8).Value = 187.9800902665885

$ws.Cells.Item(58, 7).Value = 0.03258494092351518
$ws.Cells.Item(58, 8).Value = 30.28448685524009

$ws.Cells.Item(59, 7).Value = 0.03337624549147861
$ws.Cells.Item(59, 8).Value = 40.95516268942209

$ws.Cells.Item(60, 7).Value = 0.03323841592762432
$ws.Cells.Item(60, 8).Value = 2.452813607853184

$ws.Cells.Item(61, 7).Value = 0.05423623984327209
$ws.Cells.Item(61, 8).Value = 328.4742384964268

$ws.Cells.Item(62, 7).Value = 0.04802119940148362
$ws.Cells.Item(62, 8).Value = -20.4454659788571

$ws.Cells.Item(63, 7).Value = 0.04386787898966219
$ws.Cells.Item(63, 8).Value = 34.60742861395797

$ws.Cells.Item(64, 7).Value = 0.01833743922937944
$ws.Cells.Item(64, 8).Value = -54.75162284642256

$ws.Cells.Item(65, 7).Value = 0.05243343130584844
$ws.Cells.Item(65, 8).Value = -6.472289512544452

$ws.Cells.Item(66, 7).Value = 0.08793577225891169
$ws.Cells.Item(66, 8).Value = -6.005436484931328

$ws.Cells.Item(67, 7).Value = 0.1124820930179695
$ws.Cells.Item(67, 8).Value = -2.567907552652692

$ws.Cells.Item(68, 7).Value = -0.03089695793087516
$ws.Cells.Item(68, 8).Value = 11.34424723199312

$ws.Cells.Item(69, 7).Value = -0.0215147728275465
$ws.Cells.Item(69, 8).Value = -1.380275149604262

$ws.Cells.Item(70, 7).Value = 0.09008084438607232
$ws.Cells.Item(70, 8).Value = -2.758545765491966

$ws.Cells.Item(71, 7).Value = 0.07526836716867064
$ws.Cells.Item(71, 8).Value = -17.47685505803918

$ws.Cells.Item(72, 7).Value = -0.07377374054753033
$ws.Cells.Item(72, 8).Value = -31.54848551054904

$ws.Cells.Item(73, 7).Value = -0.07309002077838081
$ws.Cells.Item(73, 8).Value = 0.9118024402693418

$ws.Cells.Item(74, 7).Value = 0.128232862303398
$ws.Cells.Item(74, 8).Value = 28.3011097756166

$ws.Cells.Item(75, 7).Value = 0.08057541477293849
$ws.Cells.Item(75, 8).Value = -17.27711797946692

$ws.Cells.Item(76, 7).Value = 0.01409206111307201
$ws.Cells.Item(76, 8).Value = -44.88959189365706

$ws.Cells.Item(77, 7).Value = 0.02303565733750857
$ws.Cells.Item(77, 8).Value = 63.26877729111258

$ws.Cells.Item(78, 7).Value = 0.08651897527265037
$ws.Cells.Item(78, 8).Value = 34.60338459365558

$ws.Cells.Item(79, 7).Value = 0.07390246641828278
$ws.Cells.Item(79, 8).Value = -3.664824704448837

$ws.Cells.Item(80, 7).Value = -0.1477189587602141
$ws.Cells.Item(80, 8).Value = 10.80225599671557

$ws.Cells.Item(81, 7).Value = -0.1633148721601642
$ws.Cells.Item(81, 8).Value = 22.26507393011026

$ws.Cells.Item(82, 7).Value = 0.1434952983243895
$ws.Cells.Item(82, 8).Value = 25.10997684239036

$ws.Cells.Item(83, 7).Value = 0.1863087717574758
$ws.Cells.Item(83, 8).Value = 4.678309777960448

$ws.Cells.Item(84, 7).Value = 0.06307673129675499
$ws.Cells.Item(84, 8).Value = 164.6096986320413

$ws.Cells.Item(85, 7).Value = 0.05678171052557069
$ws.Cells.Item(85, 8).Value = -7.785881408749473
